# "Tela de cadastro de empresa criada"
#
# The Sprint column (B) on the backlog sheet still shows the old
# "Sprint 02" label on every requirement row even though the header (B1)
# already reads just "Sprint". Bring the data rows back in line with the
# header: same text, same header formatting (bold white on blue) - i.e.
# fill B1 down over the requirement rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Find the last populated requirement row (column A), starting below the
# header row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$fillRange = $ws.Range("B2:B" + $lastRow)

# Value first ...
$fillRange.Value2 = $ws.Range("B1").Value2

# ... then bring over the header's look (font/fill/alignment) so the
# "Sprint" column matches the header style all the way down.
$ws.Range("B1").Copy()
$fillRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the author left it.
$null = $ws.Range("I16").Select()
